$p = $ppt.ActivePresentation

# 1) Notes master "datetimeFigureOut" date placeholder: 18/02/2025 -> 25/03/2025
$hf = $p.NotesMaster.HeadersFooters
$hf.DateAndTime.Text = "25/03/2025"

# 2) Slide 7 "TextBox 1" instructions box: re-created (id 2 -> 3, name "TextBox 1" -> "TextBox 2")
#    Recreate it via Copy/Paste so all paragraph/run formatting is preserved exactly,
#    then delete the original and rename the copy.
$s = $p.Slides.Item(7)
$old = $s.Shapes.Item(1)
$origHeight = $old.Height
$old.Copy()
$new = $s.Shapes.Paste()
$old.Delete()
$new.Name = "TextBox 2"

# 3) Fix typo "loose" -> "lose" in the last bullet, keeping it as a single run.
$tr = $new.TextFrame.TextRange
$full = $tr.Text
$needle = "if you RESPOND, you will loose the temporary points"
$idx = $full.IndexOf($needle)
$sub = $tr.Characters($idx + 1, $needle.Length)
$sub.Text = "if you RESPOND, you will lose the temporary points"

# The text edit re-triggers auto-fit layout on this spAutoFit textbox; restore
# the shape's original height (unchanged in the source diff). Left/Top/Width
# are untouched by the text edit so they are left alone.
$new.Height = $origHeight
